# Testing Sprint4 + creato script esecuzione
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WORK PLAN")

# Sprint 4 (rows 44-53) marked as concluded: set "CONCLUSO" label in the merged G44:G53 cell
$ws.Range("G44").Value = "CONCLUSO"

# Mark each Sprint 4 sub-task as done ("FATTO" column, checkmark)
$ws.Range("E45").Value = "✔"
$ws.Range("E46").Value = "✔"
$ws.Range("E47").Value = "✔"
$ws.Range("E48").Value = "✔"
$ws.Range("E49").Value = "✔"
$ws.Range("E50").Value = "✔"
$ws.Range("E51").Value = "✔"
$ws.Range("E52").Value = "✔"
$ws.Range("E53").Value = "✔"

# Fill in the missing completion dates ("QUANDO" column)
$ws.Range("F50").Value = 44029
$ws.Range("F51").Value = 44029
$ws.Range("F52").Value = 44030
$ws.Range("F53").Value = 44030

# Update the view to reflect scrolling/selection state after the edits
$ws.Range("C50").Select()
